$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("D7").Value = "2016-30-11 09:30:19"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("E7").Value = "2016-03-11 09:30:16"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("E7").Value = "2016-03-11 09:30:19"
